$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Row 3 is the template row (No/NoPin/NoRek/AccountSts/AccountName).
#    Clone its formatting down through row 16 before touching the
#    columns, so every new data row picks up the same per-column
#    styling that Excel would carry forward on a fill-down.
# ---------------------------------------------------------------------
for ($r = 4; $r -le 16; $r++) {
    $ws.Rows.Item(3).Copy()
    $ws.Rows.Item($r).Insert(-4121)   # xlShiftDown, carries row 3's format+values
}
$excel.CutCopyMode = 0
$ws.Range("A4:E16").ClearContents()   # keep the formatting, drop the cloned values

# Row 2 is visually "striped" (customFormat row) - alternate that
# banding onto every other new row in column A, same as rows 2/3 already do.
for ($r = 4; $r -le 16; $r += 2) {
    $ws.Cells.Item(2, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Drop the old "NoRek" column entirely. AccountSts/AccountName shift
#    left from D/E to C/D.
# ---------------------------------------------------------------------
$ws.Range("C1").EntireColumn.Delete()

# Row 2's old "NoRek" cell style (row-default) leaked into the shifted
# AccountSts cell; re-align it with the rest of the column (row 3's style).
$ws.Cells.Item(3, 3).Copy()
$ws.Cells.Item(2, 3).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. The old "AccountName" column (now D) is repurposed as the new
#    "APPDate (dd/mm/YYYY)" column - header only, no data beneath it.
# ---------------------------------------------------------------------
$ws.Range("D2:D16").Clear()
$ws.Range("D1").Value = "APPDate (dd/mm/YYYY)"

# ---------------------------------------------------------------------
# 4. Re-key every data row with the refreshed account list.
# ---------------------------------------------------------------------
$data = @(
    @(1,  "0010008956", "ACTIVE"),
    @(2,  "0320002029", "ACTIVE"),
    @(3,  "0320002041", "CLOSED - REGULER"),
    @(4,  "0320002052", "ACTIVE"),
    @(5,  "0450002380", "ACTIVE"),
    @(6,  "0500006437", "ACTIVE"),
    @(7,  "0030006614", "ACTIVE"),
    @(8,  "0400004909", "ACTIVE"),
    @(9,  "0060012241", "CLOSED - REGULER"),
    @(10, "0060012564", "ACTIVE"),
    @(11, "0450002488", "ACTIVE"),
    @(12, "0700001567", "ACTIVE"),
    @(13, "0050012007", "ACTIVE"),
    @(14, "0050012187", "ACTIVE"),
    @(15, "0030006753", "ACTIVE")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# ---------------------------------------------------------------------
# 5. Duplicate-value conditional formatting used to also watch the
#    (now deleted) NoRek column; only the NoPin rule remains, stretched
#    to cover every data row.
# ---------------------------------------------------------------------
$ws.Range("C2").FormatConditions.Delete()
$ws.Range("B2:B3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B2:B16"))

# ---------------------------------------------------------------------
# 6. Column width tweaks that came along with the edit, and leave the
#    selection where the author left it.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 9.54296875
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 21.26953125

$ws.Range("E9").Select()
